$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so they keep their
# exact textual formatting (trailing zeros, grouping dots) instead of being
# auto-converted to a Number by Excel type inference.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.050.58"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "2.416.09"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "555.30"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").Value = "142.77"
$ws.Range("E6").Value = "  +4.62%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "2.415.43"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "26.16"
$ws.Range("E14").Value = "  +6.13%  "
$ws.Range("E15").Value = "  +7.68%  "
$ws.Range("D16").Value = "2.858.82"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "62.021.63"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "2.418.76"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "323.55"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "64.91"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +5.51%  "
$ws.Range("D26").Value = "9.09"
$ws.Range("E26").Value = "  +7.43%  "
$ws.Range("D27").Value = "578.04"
$ws.Range("E27").Value = "  +17.17%  "
$ws.Range("D28").Value = "2.535.87"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  +4.41%  "
$ws.Range("E31").Value = "  +7.64%  "
$ws.Range("E32").Value = "  +6.20%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("E38").Value = "  +8.29%  "
$ws.Range("D39").Value = "0.383"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.74"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").Value = "148.31"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "41.74"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "150.92"
$ws.Range("E45").Value = "  +6.25%  "
$ws.Range("E46").Value = "  +12.16%  "
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "0.0544"
$ws.Range("E48").Value = "  +5.52%  "
$ws.Range("D49").Value = "20.31"
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("D51").Value = "0.0917"
$ws.Range("E51").Value = "  +1.85%  "
